# Update "paises" (countries) data sheet + provincias de España figures.
# The underlying data table (A3:H205) is kept sorted by "Casos totales"
# (column B) descending. This refresh:
#   - bumps a handful of countries' organic totals in place (USA, Austria,
#     Noruega, Islandia, Trinidad y Tobago)
#   - re-ranks Ucrania (jumps above Lituania/Libano/Armenia)
#   - re-ranks Camerun (brand new entrant just above Brunei) and Nigeria
#     (climbs above Honduras), pushing the countries between their old and
#     new rank down by one row
#   - refreshes the "actualizado a" timestamp banner

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 23:20"

# In-place organic updates
Set-Row 4   "Estados Unidos" 139262 15684 4435 132382 2948 225 2445
Set-Row 17  "Austria"          8774   503  479   8209  187  18   86
Set-Row 20  "Noruega"          4271   256    7   4239   91   2   25
Set-Row 45  "Islandia"         1020    57  135    883   25   0    2

# Ucrania overtakes Lituania / Libano / Armenia (rows 67-70)
Set-Row 67  "Ucrania"           475   119    6    459    0   1   10
Set-Row 68  "Lituania"          460    66    1    452    2   0    7
Set-Row 69  "Libano"            438    26   30    398    4   2   10
Set-Row 70  "Armenia"           424    17   30    391    6   2    3

# Camerun (new entrant) and Nigeria climb the ranking (rows 101-113)
Set-Row 101 "Camerun"           139    48    5    128    0   4    6
Set-Row 102 "Brunei"            126     6   34     91    1   0    1
Set-Row 103 "Afganistan"        120    10    2    114    0   0    4
Set-Row 104 "Venezuela"         119     0   39     78    2   0    2
Set-Row 105 "Sri Lanka"         117     4   11    105    5   0    1
Set-Row 106 "Nigeria"           111    14    3    107    0   0    1
Set-Row 107 "Honduras"          110    15    3    105    4   1    2
Set-Row 108 "Estado de Palestina" 108   4   18     89    0   0    1
Set-Row 109 "Mauricio"          107     5    0    104    1   1    3
Set-Row 110 "Camboya"           103     4   21     82    1   0    0
Set-Row 111 "Guadalupe"         102     0   17     83    4   0    2
Set-Row 112 "Bielorrusia"        94     0   32     62    2   0    0
Set-Row 113 "Martinica"          93     0    0     92   12   0    1

# Trinidad y Tobago organic update
Set-Row 118 "Trinidad yTobago"   78     2    1     74    0   0    3
